# Controle de Horas.xlsx -- add new "Folha Ponto" entries for 2022-12-14 (serial 44909)
# covering stages/activities that were missing, per commit message
# "More stages (all old ones are included now)".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha Ponto")

# --- 1. Seed formatting for the new rows (36-40) by copying the plain
#        Data/Inicio/Fim/Total number-format pattern from row 33, which is
#        already styled exactly like we need (date + two time cells).
$ws.Range("B33:E33").Copy()
$ws.Range("B36:E36").PasteSpecial(-4122)
$ws.Range("B37:E37").PasteSpecial(-4122)
$ws.Range("B38:E38").PasteSpecial(-4122)
$ws.Range("B39:E39").PasteSpecial(-4122)
$ws.Range("B40:E40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Row 36: SITS, 07:00 - 09:00, "Implementando fases novas"
$ws.Range("B36").Value = 44909
$ws.Range("C36").Value = 0.29166666666666669
$ws.Range("D36").Value = 0.375
$ws.Range("F36").Value = "SITS"
$ws.Range("G36").Value = "Implementando fases novas"

# --- 3. Row 37: Outros, 09:00 - 12:00, "Brainstorm de nomes pra empresa e jogos"
$ws.Range("B37").Value = 44909
$ws.Range("C37").Value = 0.375
$ws.Range("D37").Value = 0.5
$ws.Range("F37").Value = "Outros"
$ws.Range("G37").Value = "Brainstorm de nomes pra empresa e jogos"

# --- 4. Row 38: SITS, 12:00 - 12:45, "Implementando fases novas"
$ws.Range("B38").Value = 44909
$ws.Range("C38").Value = 0.5
$ws.Range("D38").Value = 0.53125
$ws.Range("F38").Value = "SITS"
$ws.Range("G38").Value = "Implementando fases novas"

# --- 5. Row 39: SITS, 13:15 - 16:00, "Implementando fases novas"
$ws.Range("B39").Value = 44909
$ws.Range("C39").Value = 0.55208333333333337
$ws.Range("D39").Value = 0.66666666666666663
$ws.Range("F39").Value = "SITS"
$ws.Range("G39").Value = "Implementando fases novas"

# --- 6. Row 40: SITS, 17:00 - 17:15, "Organizando tarefas"
$ws.Range("B40").Value = 44909
$ws.Range("C40").Value = 0.70833333333333337
$ws.Range("D40").Value = 0.71875
$ws.Range("F40").Value = "SITS"
$ws.Range("G40").Value = "Organizando tarefas"

# --- 7. Restore the previously-selected cell pointer (now D41, the next
#        empty row) like the author's Excel session left it.
$ws.Range("D41").Select()

$excel.Calculate()
